$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$values = @(0, 1, 0, 0, 0, 0, 0, 0, 0, 0)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $values[$i]
}
